$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 6515.857  # H18 6515.7856 -> 6515.857
$ws.Cells.Item(18, 9).Value = 3851.8333  # I18 3851.75 -> 3851.8333
$ws.Cells.Item(18, 11).Value = 3851.8333  # K18 3851.75 -> 3851.8333
$ws.Cells.Item(18, 13).Value = -3567.8333  # M18 -3567.75 -> -3567.8333

$ws.Cells.Item(19, 8).Value = 1433.3  # H19 1587.0555 -> 1433.3
$ws.Cells.Item(19, 9).Value = 1481.2727  # I19 1559.5 -> 1481.2727
$ws.Cells.Item(19, 10).Value = 1374.6666  # J19 1621.5 -> 1374.6666
$ws.Cells.Item(19, 11).Value = 1481.2727  # K19 1559.5 -> 1481.2727
$ws.Cells.Item(19, 12).Value = 1374.6666  # L19 1621.5 -> 1374.6666
$ws.Cells.Item(19, 13).Value = -1306.2727  # M19 -1384.5 -> -1306.2727
$ws.Cells.Item(19, 14).Value = -1724.6666  # N19 -1971.5 -> -1724.6666

$ws.Cells.Item(40, 8).Value = 2071.5715  # H40 1416.8334 -> 2071.5715
$ws.Cells.Item(40, 10).Value = 3500  # J40 1000 -> 3500
$ws.Cells.Item(40, 12).Value = 3500  # L40 1000 -> 3500
$ws.Cells.Item(40, 14).Value = -3850  # N40 -1350 -> -3850

$ws.Cells.Item(64, 8).Value = 28598.666  # H64 26738.7 -> 28598.666
$ws.Cells.Item(64, 10).Value = 10000  # J64 9999.5 -> 10000
$ws.Cells.Item(64, 12).Value = 10000  # L64 9999.5 -> 10000
$ws.Cells.Item(64, 14).Value = -10496  # N64 -10495.5 -> -10496

$ws.Cells.Item(67, 8).Value = 28598.666  # H67 26738.7 -> 28598.666
$ws.Cells.Item(67, 10).Value = 10000  # J67 9999.5 -> 10000
$ws.Cells.Item(67, 12).Value = 10000  # L67 9999.5 -> 10000
$ws.Cells.Item(67, 14).Value = -11716  # N67 -11715.5 -> -11716

$ws.Cells.Item(86, 8).Value = 2575.1667  # H86 2850 -> 2575.1667
$ws.Cells.Item(86, 9).Value = 1860.7142  # I86 1943.25 -> 1860.7142
$ws.Cells.Item(86, 11).Value = 1860.7142  # K86 1943.25 -> 1860.7142
$ws.Cells.Item(86, 13).Value = -737.7141999999999  # M86 -820.25 -> -737.7141999999999

$ws.Cells.Item(89, 8).Value = 2575.1667  # H89 2850 -> 2575.1667
$ws.Cells.Item(89, 9).Value = 1860.7142  # I89 1943.25 -> 1860.7142
$ws.Cells.Item(89, 11).Value = 9303.571  # K89 9716.25 -> 9303.571
$ws.Cells.Item(89, 13).Value = -3687.571  # M89 -4100.25 -> -3687.571

$ws.Cells.Item(98, 8).Value = 3845.3333  # H98 4056.353 -> 3845.3333
$ws.Cells.Item(98, 9).Value = 1314.3448  # I98 1392.5927 -> 1314.3448
$ws.Cells.Item(98, 11).Value = 1314.3448  # K98 1392.5927 -> 1314.3448
$ws.Cells.Item(98, 13).Value = 183.6551999999999  # M98 105.4073000000001 -> 183.6551999999999

$ws.Cells.Item(122, 8).Value = 3845.3333  # H122 4056.353 -> 3845.3333
$ws.Cells.Item(122, 9).Value = 1314.3448  # I122 1392.5927 -> 1314.3448
$ws.Cells.Item(122, 11).Value = 3943.0344  # K122 4177.7781 -> 3943.0344
$ws.Cells.Item(122, 13).Value = -1493.0344  # M122 -1727.7781 -> -1493.0344

$ws.Cells.Item(125, 8).Value = 40216  # H125 45832.855 -> 40216
$ws.Cells.Item(125, 10).Value = 2139.4  # J125 2449.75 -> 2139.4
$ws.Cells.Item(125, 12).Value = 19254.6  # L125 22047.75 -> 19254.6
$ws.Cells.Item(125, 14).Value = -24174.6  # N125 -26967.75 -> -24174.6

$ws.Cells.Item(129, 8).Value = 1429.75  # H129 1200.25 -> 1429.75
$ws.Cells.Item(129, 9).Value = 832.44446  # I129 773.7 -> 832.44446
$ws.Cells.Item(129, 10).Value = 3221.6667  # J129 3333 -> 3221.6667
$ws.Cells.Item(129, 11).Value = 2497.33338  # K129 2321.1 -> 2497.33338
$ws.Cells.Item(129, 12).Value = 9665.000100000001  # L129 9999 -> 9665.000100000001
$ws.Cells.Item(129, 13).Value = 2502.66662  # M129 2678.9 -> 2502.66662
$ws.Cells.Item(129, 14).Value = -19665.0001  # N129 -19999 -> -19665.0001

$ws.Cells.Item(132, 8).Value = 25123.066  # H132 25964.586 -> 25123.066
$ws.Cells.Item(132, 9).Value = 26831.857  # I132 27799 -> 26831.857
$ws.Cells.Item(132, 11).Value = 80495.571  # K132 83397 -> 80495.571
$ws.Cells.Item(132, 13).Value = -77965.571  # M132 -80867 -> -77965.571

$ws.Cells.Item(133, 8).Value = 89780  # H133 0 -> 89780
$ws.Cells.Item(133, 10).Value = 89780  # J133 0 -> 89780
$ws.Cells.Item(133, 12).Value = 89780  # L133 0 -> 89780
$ws.Cells.Item(133, 14).Value = -99900  # N133 None -> -99900

$ws.Cells.Item(137, 8).Value = 101617.336  # H137 68301.89 -> 101617.336
$ws.Cells.Item(137, 9).Value = 100634  # I137 60989.2 -> 100634
$ws.Cells.Item(137, 10).Value = 102600.664  # J137 77442.75 -> 102600.664
$ws.Cells.Item(137, 11).Value = 301902  # K137 182967.6 -> 301902
$ws.Cells.Item(137, 12).Value = 307801.992  # L137 232328.25 -> 307801.992
$ws.Cells.Item(137, 13).Value = -299352  # M137 -180417.6 -> -299352
$ws.Cells.Item(137, 14).Value = -312901.992  # N137 -237428.25 -> -312901.992

$ws.Cells.Item(138, 8).Value = 17676.863  # H138 17939.584 -> 17676.863
$ws.Cells.Item(138, 9).Value = 1648.862  # I138 1686.3214 -> 1648.862
$ws.Cells.Item(138, 11).Value = 4946.586  # K138 5058.9642 -> 4946.586
$ws.Cells.Item(138, 13).Value = 193.4139999999998  # M138 81.03579999999965 -> 193.4139999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2892.842  # H45 2578.7083 -> 2892.842
$ws.Cells.Item(45, 9).Value = 2360.4546  # I45 2150.9285 -> 2360.4546
$ws.Cells.Item(45, 10).Value = 3624.875  # J45 3177.6 -> 3624.875
$ws.Cells.Item(45, 11).Value = 2360.4546  # K45 2150.9285 -> 2360.4546
$ws.Cells.Item(45, 12).Value = 3624.875  # L45 3177.6 -> 3624.875
$ws.Cells.Item(45, 13).Value = -1983.4546  # M45 -1773.9285 -> -1983.4546
$ws.Cells.Item(45, 14).Value = -4378.875  # N45 -3931.6 -> -4378.875

$ws.Cells.Item(63, 8).Value = 4506.125  # H63 4581.0625 -> 4506.125
$ws.Cells.Item(63, 9).Value = 4016.3333  # I63 4459.6 -> 4016.3333
$ws.Cells.Item(63, 10).Value = 4800  # J63 4636.273 -> 4800
$ws.Cells.Item(63, 11).Value = 4016.3333  # K63 4459.6 -> 4016.3333
$ws.Cells.Item(63, 12).Value = 4800  # L63 4636.273 -> 4800
$ws.Cells.Item(63, 13).Value = -3330.3333  # M63 -3773.6 -> -3330.3333
$ws.Cells.Item(63, 14).Value = -6172  # N63 -6008.273 -> -6172

$ws.Cells.Item(66, 8).Value = 4506.125  # H66 4581.0625 -> 4506.125
$ws.Cells.Item(66, 9).Value = 4016.3333  # I66 4459.6 -> 4016.3333
$ws.Cells.Item(66, 10).Value = 4800  # J66 4636.273 -> 4800
$ws.Cells.Item(66, 11).Value = 20081.6665  # K66 22298 -> 20081.6665
$ws.Cells.Item(66, 12).Value = 24000  # L66 23181.365 -> 24000
$ws.Cells.Item(66, 13).Value = -16649.6665  # M66 -18866 -> -16649.6665
$ws.Cells.Item(66, 14).Value = -30864  # N66 -30045.365 -> -30864

$ws.Cells.Item(122, 8).Value = 1928.7916  # H122 2175.2632 -> 1928.7916
$ws.Cells.Item(122, 9).Value = 1662.1052  # I122 1901.3572 -> 1662.1052
$ws.Cells.Item(122, 11).Value = 4986.3156  # K122 5704.071599999999 -> 4986.3156
$ws.Cells.Item(122, 13).Value = -2536.3156  # M122 -3254.071599999999 -> -2536.3156

$ws.Cells.Item(132, 8).Value = 1562.4615  # H132 1601 -> 1562.4615
$ws.Cells.Item(132, 9).Value = 1328.3636  # I132 1363.0952 -> 1328.3636
$ws.Cells.Item(132, 11).Value = 3985.0908  # K132 4089.2856 -> 3985.0908
$ws.Cells.Item(132, 13).Value = -1455.0908  # M132 -1559.2856 -> -1455.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 10433.591  # H20 10890.81 -> 10433.591
$ws.Cells.Item(20, 9).Value = 12306  # I20 12980.941 -> 12306
$ws.Cells.Item(20, 11).Value = 12306  # K20 12980.941 -> 12306
$ws.Cells.Item(20, 13).Value = -12059  # M20 -12733.941 -> -12059

$ws.Cells.Item(22, 8).Value = 922.3333  # H22 917.1429000000001 -> 922.3333
$ws.Cells.Item(22, 10).Value = 0  # J22 886 -> 0
$ws.Cells.Item(22, 12).Value = 0  # L22 886 -> 0
$ws.Cells.Item(22, 14).ClearContents()  # N22 remove (was -1232)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 1686.8636  # H94 1689.3636 -> 1686.8636
$ws.Cells.Item(94, 9).Value = 1384.5555  # I94 1385.7778 -> 1384.5555
$ws.Cells.Item(94, 10).Value = 1896.1538  # J94 1899.5385 -> 1896.1538
$ws.Cells.Item(94, 11).Value = 1384.5555  # K94 1385.7778 -> 1384.5555
$ws.Cells.Item(94, 12).Value = 1896.1538  # L94 1899.5385 -> 1896.1538
$ws.Cells.Item(94, 13).Value = -933.5554999999999  # M94 -934.7778000000001 -> -933.5554999999999
$ws.Cells.Item(94, 14).Value = -2798.1538  # N94 -2801.5385 -> -2798.1538

$ws.Cells.Item(105, 8).Value = 1483  # H105 2193.25 -> 1483
$ws.Cells.Item(105, 9).Value = 1043.375  # I105 1257.8334 -> 1043.375
$ws.Cells.Item(105, 10).Value = 5000  # J105 4999.5 -> 5000
$ws.Cells.Item(105, 11).Value = 1043.375  # K105 1257.8334 -> 1043.375
$ws.Cells.Item(105, 12).Value = 5000  # L105 4999.5 -> 5000
$ws.Cells.Item(105, 13).Value = 703.625  # M105 489.1666 -> 703.625
$ws.Cells.Item(105, 14).Value = -8494  # N105 -8493.5 -> -8494

$ws.Cells.Item(134, 8).Value = 2385  # H134 2457.7083 -> 2385
$ws.Cells.Item(134, 9).Value = 1636.5  # I134 1687.1765 -> 1636.5
$ws.Cells.Item(134, 10).Value = 4069.125  # J134 4329 -> 4069.125
$ws.Cells.Item(134, 11).Value = 4909.5  # K134 5061.529500000001 -> 4909.5
$ws.Cells.Item(134, 12).Value = 12207.375  # L134 12987 -> 12207.375
$ws.Cells.Item(134, 13).Value = -2374.5  # M134 -2526.529500000001 -> -2374.5
$ws.Cells.Item(134, 14).Value = -17277.375  # N134 -18057 -> -17277.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 746.44446  # H5 765.4286 -> 746.44446
$ws.Cells.Item(5, 9).Value = 715  # I5 750 -> 715
$ws.Cells.Item(5, 11).Value = 2145  # K5 2250 -> 2145
$ws.Cells.Item(5, 13).Value = -2033  # M5 -2138 -> -2033

$ws.Cells.Item(12, 8).Value = 180.27272  # H12 196.95 -> 180.27272
$ws.Cells.Item(12, 9).Value = 205.33333  # I12 245.4 -> 205.33333
$ws.Cells.Item(12, 10).Value = 170.875  # J12 180.8 -> 170.875
$ws.Cells.Item(12, 11).Value = 615.99999  # K12 736.2 -> 615.99999
$ws.Cells.Item(12, 12).Value = 512.625  # L12 542.4000000000001 -> 512.625
$ws.Cells.Item(12, 13).Value = -442.99999  # M12 -563.2 -> -442.99999
$ws.Cells.Item(12, 14).Value = -858.625  # N12 -888.4000000000001 -> -858.625

$ws.Cells.Item(135, 8).Value = 746.44446  # H135 765.4286 -> 746.44446
$ws.Cells.Item(135, 9).Value = 715  # I135 750 -> 715
$ws.Cells.Item(135, 11).Value = 6435  # K135 6750 -> 6435
$ws.Cells.Item(135, 13).Value = -3900  # M135 -4215 -> -3900

$ws.Cells.Item(137, 8).Value = 3277.5715  # H137 3491.625 -> 3277.5715
$ws.Cells.Item(137, 10).Value = 3933.25  # J137 4144.6 -> 3933.25
$ws.Cells.Item(137, 12).Value = 11799.75  # L137 12433.8 -> 11799.75
$ws.Cells.Item(137, 14).Value = -21999.75  # N137 -22633.8 -> -21999.75

$ws.Cells.Item(140, 8).Value = 2166.4285  # H140 2464.1667 -> 2166.4285
$ws.Cells.Item(140, 9).Value = 2166.4285  # I140 2464.1667 -> 2166.4285
$ws.Cells.Item(140, 11).Value = 6499.2855  # K140 7392.500100000001 -> 6499.2855
$ws.Cells.Item(140, 13).Value = -1319.2855  # M140 -2212.500100000001 -> -1319.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6916.2856  # H70 6614.25 -> 6916.2856
$ws.Cells.Item(70, 10).Value = 7139.3335  # J70 6479.5 -> 7139.3335
$ws.Cells.Item(70, 12).Value = 7139.3335  # L70 6479.5 -> 7139.3335
$ws.Cells.Item(70, 14).Value = -7679.3335  # N70 -7019.5 -> -7679.3335

$ws.Cells.Item(73, 8).Value = 6916.2856  # H73 6614.25 -> 6916.2856
$ws.Cells.Item(73, 10).Value = 7139.3335  # J73 6479.5 -> 7139.3335
$ws.Cells.Item(73, 12).Value = 7139.3335  # L73 6479.5 -> 7139.3335
$ws.Cells.Item(73, 14).Value = -9011.333500000001  # N73 -8351.5 -> -9011.333500000001

$ws.Cells.Item(80, 8).Value = 5461.722  # H80 5200.5557 -> 5461.722
$ws.Cells.Item(80, 9).Value = 1837.5  # I80 1766.5555 -> 1837.5
$ws.Cells.Item(80, 10).Value = 8361.1  # J80 8634.556 -> 8361.1
$ws.Cells.Item(80, 11).Value = 1837.5  # K80 1766.5555 -> 1837.5
$ws.Cells.Item(80, 12).Value = 8361.1  # L80 8634.556 -> 8361.1
$ws.Cells.Item(80, 13).Value = -839.5  # M80 -768.5554999999999 -> -839.5
$ws.Cells.Item(80, 14).Value = -10357.1  # N80 -10630.556 -> -10357.1

$ws.Cells.Item(83, 8).Value = 5461.722  # H83 5200.5557 -> 5461.722
$ws.Cells.Item(83, 9).Value = 1837.5  # I83 1766.5555 -> 1837.5
$ws.Cells.Item(83, 10).Value = 8361.1  # J83 8634.556 -> 8361.1
$ws.Cells.Item(83, 11).Value = 9187.5  # K83 8832.7775 -> 9187.5
$ws.Cells.Item(83, 12).Value = 41805.5  # L83 43172.78 -> 41805.5
$ws.Cells.Item(83, 13).Value = -4195.5  # M83 -3840.7775 -> -4195.5
$ws.Cells.Item(83, 14).Value = -51789.5  # N83 -53156.78 -> -51789.5

$ws.Cells.Item(102, 8).Value = 2091.3684  # H102 2117.875 -> 2091.3684
$ws.Cells.Item(102, 9).Value = 1267  # I102 1324 -> 1267
$ws.Cells.Item(102, 10).Value = 4399.6  # J102 4499.5 -> 4399.6
$ws.Cells.Item(102, 11).Value = 1267  # K102 1324 -> 1267
$ws.Cells.Item(102, 12).Value = 4399.6  # L102 4499.5 -> 4399.6
$ws.Cells.Item(102, 13).Value = 355  # M102 298 -> 355
$ws.Cells.Item(102, 14).Value = -7643.6  # N102 -7743.5 -> -7643.6

$ws.Cells.Item(126, 8).Value = 3634.2856  # H126 3018.5557 -> 3634.2856
$ws.Cells.Item(126, 9).Value = 3088  # I126 2452.4285 -> 3088
$ws.Cells.Item(126, 11).Value = 9264  # K126 7357.2855 -> 9264
$ws.Cells.Item(126, 13).Value = -6794  # M126 -4887.2855 -> -6794

$ws.Cells.Item(132, 8).Value = 2671.5  # H132 2639.4285 -> 2671.5
$ws.Cells.Item(132, 9).Value = 2377.5454  # I132 2446.7144 -> 2377.5454
$ws.Cells.Item(132, 10).Value = 3210.4167  # J132 2928.5 -> 3210.4167
$ws.Cells.Item(132, 11).Value = 7132.6362  # K132 7340.1432 -> 7132.6362
$ws.Cells.Item(132, 12).Value = 9631.250100000001  # L132 8785.5 -> 9631.250100000001
$ws.Cells.Item(132, 13).Value = -4602.6362  # M132 -4810.1432 -> -4602.6362
$ws.Cells.Item(132, 14).Value = -14691.2501  # N132 -13845.5 -> -14691.2501

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1004  # H22 1005.46875 -> 1004
$ws.Cells.Item(22, 10).Value = 1057.6316  # J22 1060.1052 -> 1057.6316
$ws.Cells.Item(22, 12).Value = 1057.6316  # L22 1060.1052 -> 1057.6316
$ws.Cells.Item(22, 14).Value = -1647.6316  # N22 -1650.1052 -> -1647.6316

$ws.Cells.Item(27, 8).Value = 1004  # H27 1005.46875 -> 1004
$ws.Cells.Item(27, 10).Value = 1057.6316  # J27 1060.1052 -> 1057.6316
$ws.Cells.Item(27, 12).Value = 1057.6316  # L27 1060.1052 -> 1057.6316
$ws.Cells.Item(27, 14).Value = -1271.6316  # N27 -1274.1052 -> -1271.6316

$ws.Cells.Item(133, 8).Value = 105163  # H133 0 -> 105163
$ws.Cells.Item(133, 10).Value = 105163  # J133 0 -> 105163
$ws.Cells.Item(133, 12).Value = 105163  # L133 0 -> 105163
$ws.Cells.Item(133, 14).Value = -110223  # N133 None -> -110223

$ws.Cells.Item(136, 8).Value = 4823.278  # H136 4813.8237 -> 4823.278
$ws.Cells.Item(136, 10).Value = 5708.4287  # J136 5829.1665 -> 5708.4287
$ws.Cells.Item(136, 12).Value = 17125.2861  # L136 17487.4995 -> 17125.2861
$ws.Cells.Item(136, 14).Value = -22225.2861  # N136 -22587.4995 -> -22225.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(52, 8).Value = 0  # H52 49999 -> 0
$ws.Cells.Item(52, 10).Value = 0  # J52 49999 -> 0
$ws.Cells.Item(52, 12).Value = 0  # L52 49999 -> 0
$ws.Cells.Item(52, 14).ClearContents()  # N52 remove (was -50451)

$ws.Cells.Item(81, 8).Value = 8774  # H81 8928.454 -> 8774
$ws.Cells.Item(81, 9).Value = 9301.4  # I81 10134.889 -> 9301.4
$ws.Cells.Item(81, 10).Value = 3500  # J81 3499.5 -> 3500
$ws.Cells.Item(81, 11).Value = 18602.8  # K81 20269.778 -> 18602.8
$ws.Cells.Item(81, 12).Value = 7000  # L81 6999 -> 7000
$ws.Cells.Item(81, 13).Value = -17541.8  # M81 -19208.778 -> -17541.8
$ws.Cells.Item(81, 14).Value = -9122  # N81 -9121 -> -9122

$ws.Cells.Item(84, 8).Value = 8774  # H84 8928.454 -> 8774
$ws.Cells.Item(84, 9).Value = 9301.4  # I84 10134.889 -> 9301.4
$ws.Cells.Item(84, 10).Value = 3500  # J84 3499.5 -> 3500
$ws.Cells.Item(84, 11).Value = 93014  # K84 101348.89 -> 93014
$ws.Cells.Item(84, 12).Value = 35000  # L84 34995 -> 35000
$ws.Cells.Item(84, 13).Value = -87710  # M84 -96044.88999999998 -> -87710
$ws.Cells.Item(84, 14).Value = -45608  # N84 -45603 -> -45608
